$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.760.25"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.48%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.106.48"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.22"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.69"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.106.59"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.93%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.521"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.158"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.50%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.16"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.48%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.491"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +4.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000236"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.54"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.32%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.615.00"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.72%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.720.60"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.100.00"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.35%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.95"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "497.05"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.55"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +3.89%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.18"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +11.66%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.705"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +3.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.65"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.91%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.45"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.16%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.86"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.44"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +3.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.13"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "27.24"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.73%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.996"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.65"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +6.03%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.88%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.95"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +6.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.43"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.84%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.35"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.12%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0419"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.94%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "451.97"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.78%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0837"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.03%  "
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.86"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -5.62%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.073.10"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +3.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.46"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.83%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.71%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.284"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +8.83%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "28.69"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.31"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +7.08%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.115"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₃0537"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.76%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.20"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +5.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "118.30"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.00%  "
